# "added ack for debugger"
# Fills in the newly-documented FPGA debug pipeline signals (0x18-0x1D range)
# on the DebugDataTransfer sheet, and leaves the selection on B16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 - 0x18 pipeline_immediate_out (16 Bit)
$ws.Range("B26").Value = "0x18"
$ws.Range("E26").Value = "pipeline_immediate_out"
$ws.Range("G26").Value = "0x18"
$ws.Range("H26").Value = "16 Bit"
$ws.Range("J26").Value = "pipeline_immediate_out"

# Row 27 - 0x19 pipeline_write_address_out (4 Bit)
$ws.Range("B27").Value = "0x19"
$ws.Range("E27").Value = "pipeline_write_address_out"
$ws.Range("G27").Value = "0x19"
$ws.Range("H27").Value = "4 Bit"
$ws.Range("J27").Value = "pipeline_write_address_out"

# Row 28 - 0x1A pipeline_whb_wlb_out (2 Bit)
$ws.Range("B28").Value = "0x1A"
$ws.Range("E28").Value = "pipeline_whb_wlb_out"
$ws.Range("G28").Value = "0x1A"
$ws.Range("H28").Value = "2 Bit"
$ws.Range("J28").Value = "pipeline_whb_wlb_out"

# Rows 29 (0x1B) and 30 (0x1C) already contain their data - untouched.

# Row 31 - 0x1D pipeline_is_alu_ram_gpu_op_out (3 Bit)
$ws.Range("B31").Value = "0x1D"
$ws.Range("E31").Value = "pipeline_is_alu_ram_gpu_op_out"
$ws.Range("G31").Value = "0x1D"
$ws.Range("H31").Value = "3 Bit"
$ws.Range("J31").Value = "pipeline_is_alu_ram_gpu_op_out"

# Update the view: scroll back to top and select B16 (matches the author's
# final cursor position while reviewing the new ack/debug rows).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B16").Select()
